$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header / masthead text updates ---
$ws.Range("M6").Value = "Jessica S. Tisch"
$ws.Range("A8").Value = "Volume 31   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -42.857142857142
$ws.Range("L15").Value = -11.111111111111
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 75
$ws.Range("I16").Value = 114
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 14
$ws.Range("L16").Value = 9.615384615384
$ws.Range("M16").Value = -43.842364532019
$ws.Range("N16").Value = -82.298136645962
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -61.538461538461
$ws.Range("I17").Value = 159
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 6
$ws.Range("L17").Value = 1.273885350318
$ws.Range("M17").Value = 32.5
$ws.Range("N17").Value = -46.464646464646
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 69
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = -17.857142857142
$ws.Range("L18").Value = -37.272727272727
$ws.Range("M18").Value = -73.461538461538
$ws.Range("N18").Value = -93.496701225259
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -8.510638297872
$ws.Range("I19").Value = 531
$ws.Range("J19").Value = 607
$ws.Range("K19").Value = -12.520593080724
$ws.Range("L19").Value = 9.034907597535
$ws.Range("M19").Value = 12.025316455696
$ws.Range("N19").Value = -10.906040268456
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 57.142857142857
$ws.Range("I20").Value = 161
$ws.Range("J20").Value = 122
$ws.Range("K20").Value = 31.967213114754
$ws.Range("L20").Value = 31.967213114754
$ws.Range("M20").Value = 8.783783783783
$ws.Range("N20").Value = -93.633847370502
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -19.047619047619
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = -2.469135802469
$ws.Range("I21").Value = 1044
$ws.Range("J21").Value = 1080
$ws.Range("K21").Value = -3.333333333333
$ws.Range("L21").Value = 5.454545454545
$ws.Range("M21").Value = -14.636140637776
$ws.Range("N21").Value = -79.787028073572
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("L23").Value = -12.903225806451
$ws.Range("M23").Value = 3.846153846153
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 65
$ws.Range("F24").Value = 141
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = 21.551724137931
$ws.Range("I24").Value = 1524
$ws.Range("J24").Value = 1211
$ws.Range("K24").Value = 25.846407927332
$ws.Range("L24").Value = 44.592030360531
$ws.Range("M24").Value = 60.929250263991
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 127.272727272727
$ws.Range("F25").Value = 110
$ws.Range("G25").Value = 89
$ws.Range("H25").Value = 23.595505617977
$ws.Range("I25").Value = 1257
$ws.Range("J25").Value = 891
$ws.Range("K25").Value = 41.077441077441
$ws.Range("L25").Value = 88.738738738738
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 12
$ws.Range("I26").Value = 291
$ws.Range("J26").Value = 253
$ws.Range("K26").Value = 15.019762845849
$ws.Range("L26").Value = 35.348837209302
$ws.Range("M26").Value = -11.818181818181
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = -44
$ws.Range("L27").Value = -6.666666666666
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = 19.354838709677

# --- Cells switching from text placeholder ("0"/"***.*") to numeric ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Cells switching from numeric to text placeholder ("0"/"***.*") ---
# Row 22 (Transit) is an untouched reference row already styled/typed this way.
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$excel.CutCopyMode = 0
